$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are written as text (matching inlineStr source),
# not auto-converted to numbers by Excel.
$textCells = @("D5","D6","D9","D10","D11","D13","D14","D16","D20","D21","D22","D24","D28","D29","D30","D34","D35","D37","D38","D39","D40","D43","D45","D46","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.120.85"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.653.12"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "218.33"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "0.5291"
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  -2.32%  "
$ws.Range("D9").Value = "0.06328"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "20.39"
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("D11").Value = "0.07757"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "1.687.00"
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").Value = "4.491"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "0.5462"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "0.0₅8127"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").Value = "65.22"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "26.111.59"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("D20").Value = "193.36"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "10.03"
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("D22").Value = "5.991"
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "140.45"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").Value = "1.434"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").Value = "0.05937"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").Value = "1.275"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("E31").Value = "  -3.06%  "
$ws.Range("E32").Value = "  -2.46%  "
$ws.Range("E33").Value = "  -5.35%  "
$ws.Range("D34").Value = "2.413"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").Value = "0.9449"
$ws.Range("E35").Value = "  -3.56%  "
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").Value = "0.5624"
$ws.Range("E37").Value = "  -4.66%  "
$ws.Range("D38").Value = "0.01608"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").Value = "5.843"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("D40").Value = "0.8465"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.010.47"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "100.98"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").Value = "1.801.03"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "56.85"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").Value = "1.004"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  -8.68%  "
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.05150"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.469"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").Value = "7.721"
$ws.Range("E51").Value = "  -4.82%  "
